# Updated cryptos list on Mon Dec 11 14:37:43 UTC 2023 with GitHub Actions
# Refresh Price (col D) / Volume(1h) (col E) figures, and swap a couple of
# coin rows (THORChain<->MultiversX, FTXToken<->FraxShare) plus replace
# ARBITRUM with BitTorrent-New in the last row.
#
# Note: several "Price" strings look numeric (e.g. "240.73") but must stay
# plain text cells (as they were before), matching the source feed's
# formatting. Prefixing with a leading apostrophe forces Excel to store the
# value as text instead of auto-converting it to a number, and resetting
# the range Style back to "Normal" afterwards avoids leaving a stray
# quote-prefix/number-format style on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.926.35'
$ws.Range("E2").Value = '  -4.23%  '

$ws.Range("D3").Value = '2.223.08'
$ws.Range("E3").Value = '  -5.28%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").Value = "'240.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("D6").Value = "'0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.12%  '

$ws.Range("D7").Value = "'67.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -8.01%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").Value = "'0.546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.23%  '

$ws.Range("D10").Value = "'0.0961"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.47%  '

$ws.Range("D11").Value = "'58.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.72%  '

$ws.Range("D12").Value = "'35.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.47%  '

$ws.Range("E13").Value = '  -2.75%  '

$ws.Range("D14").Value = "'6.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.37%  '

$ws.Range("D15").Value = '2.560.89'
$ws.Range("E15").Value = '  -5.02%  '

$ws.Range("D16").Value = "'14.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.44%  '

$ws.Range("D17").Value = "'0.846"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.88%  '

$ws.Range("D18").Value = '2.228.15'
$ws.Range("E18").Value = '  -4.98%  '

$ws.Range("D19").Value = '41.825.13'
$ws.Range("E19").Value = '  -4.46%  '

$ws.Range("D20").Value = '0.0₃0950'
$ws.Range("E20").Value = '  -7.84%  '

$ws.Range("D21").Value = "'72.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.79%  '

$ws.Range("D22").Value = "'6.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -8.40%  '

$ws.Range("D23").Value = "'233.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.78%  '

$ws.Range("D24").Value = "'2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +12.22%  '

$ws.Range("E25").Value = '  -0.15%  '

$ws.Range("D26").Value = "'3.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.18%  '

$ws.Range("D27").Value = "'2.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.87%  '

$ws.Range("E28").Value = '  -3.42%  '

$ws.Range("D29").Value = "'9.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.90%  '

$ws.Range("D30").Value = "'170.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.54%  '

$ws.Range("D31").Value = "'20.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.64%  '

$ws.Range("E32").Value = '  -5.42%  '

$ws.Range("E33").Value = '  -6.40%  '

$ws.Range("D34").Value = "'0.0708"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.99%  '

$ws.Range("D35").Value = "'5.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.34%  '

$ws.Range("D36").Value = "'4.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.18%  '

$ws.Range("E37").Value = '  -0.45%  '

$ws.Range("D38").Value = "'23.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +22.97%  '

$ws.Range("D39").Value = "'0.0278"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.88%  '

$ws.Range("D40").Value = "'2.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.08%  '

$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").Value = "'66.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.39%  '

$ws.Range("B42").Value = 'THORChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D42").Value = "'5.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.25%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = "'9.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.55%  '

$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").Value = "'4.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -13.73%  '

$ws.Range("E45").Value = '  -4.41%  '

$ws.Range("D46").Value = "'0.188"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.11%  '

$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("D48").Value = "'4.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.59%  '

$ws.Range("E49").Value = '  -3.83%  '

$ws.Range("E50").Value = '  -2.40%  '

$ws.Range("B51").Value = 'BitTorrent-New'
$ws.Range("C51").Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range("D51").Value = '0.0₃0146'
$ws.Range("E51").Value = '  +3.71%  '
